$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "重要管脚定义"

# --- New "信号名称" (signal name) / "定义" (definition) block in G1:H4 ---

# Give G1:G4 the same bordered style already used throughout the table
# (style index 1) by copying the format from an existing bordered cell.
$ws.Range("A1").Copy()
$ws.Range("G1:G4").PasteSpecial(-4122)

# Column H: first apply wrap-text at the column level (this mints the
# "no border + wrap" style that becomes the column's default/style index),
# then paste the existing bordered style onto H1:H4 and re-apply wrap so
# those four cells end up both bordered and wrapped.
$ws.Columns.Item(8).WrapText = $true
$ws.Range("A1").Copy()
$ws.Range("H1:H4").PasteSpecial(-4122)
$ws.Range("H1:H4").WrapText = $true

# Cell values - written in an order that reproduces the canonical
# shared-string table order:
#   32=信号名称, 33=SDIO, 34=SCLK definition, 35=SDIO definition, 36=CSB definition
$ws.Range("G1").Value = "信号名称"
$ws.Range("G3").Value = "SDIO"

$ws.Range("H1").Value = "定义"
$ws.Range("H2").Value = "SPI串行时钟，用于对串行接口的读写进行同步"
$ws.Range("H3").Value = "SPI串行数据输入输出，一个复用的信号可以作为输入或输出，取决于接收到的指令和时序的位置"
$ws.Range("H4").Value = "SPI片选，低有效，可以控制读写循环"

$ws.Range("G2").Value = "SCLK"
$ws.Range("G4").Value = "CSB"

# Column H width (character-width 34.1667 stores as width 35, matching the
# target column definition exactly).
$ws.Columns.Item(8).ColumnWidth = 34.1666666666666664

# Row heights for rows 2-4 so the wrapped definitions fit.
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 30

# Selection matches the newly-added block.
$ws.Range("G1:H4").Select()
